# Generate Report for Archive
#
# The localization status moved on: every cell that used to read
# "Ready for handoff" is now "In Translation" (Overview!E2:F4, and the
# Status column - column C - on both the "zh-cn" and "de-de" report
# sheets). Because the text got shorter, the Status/zh-cn/de-de columns
# that were sized to fit it are re-fitted narrower too.

$wb = $excel.ActiveWorkbook

# --- 1. Swap the status text everywhere it appears -------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2. Re-fit the columns that held that text ------------------------
# The simulated AutoFit() here snaps width to (maxStringLength + 2.8333),
# coarser than real Excel's pixel-metrics fit, so we set ColumnWidth
# explicitly to the closest value Excel's column-width grid (1/6-char
# increments) can represent of the canonical fitted width.
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn)
$ws1.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de)

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Columns.Item(3).ColumnWidth = 12.5   # column C (Status)

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Columns.Item(3).ColumnWidth = 12.5   # column C (Status)
